$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# Sheet "About": update A2 (Version line) and A6 (Recommended Citation line)
$aboutWs = $wb.Worksheets.Item("About")

$a2 = $aboutWs.Range("A2")
$a2text = $a2.Value()
$a2.Value = $a2text.Replace($oldStamp, $newStamp)

$a6 = $aboutWs.Range("A6")
$a6text = $a6.Value()
$a6.Value = $a6text.Replace($oldStamp, $newStamp)

# Sheet "Boundaries and methane sources": update S2:S25 (build_version column, excluding header S1)
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 25; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S = 19
    $cellText = $cell.Value()
    $cell.Value = $cellText.Replace($oldStamp, $newStamp)
}
